$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changesData = @"
3|F|Benfica|text
10|F|Portimonense|text
14|B|6876465|num
14|E|Benfica|text
14|F|Estrela|text
14|G|2|num
14|H|0|num
14|J|0|num
14|K|H|text
14|L|1.125|num
14|M|8|num
14|N|19|num
14|O|1.142|num
14|P|7.5|num
14|Q|17|num
14|R|-2.25|num
14|S|1.925|num
14|T|1.925|num
14|U|3.5|num
14|V|1.875|num
14|W|1.975|num
14|X|0.1419999999999999|num
14|Z|-1|num
14|AA|-0.5|num
14|AB|0.4625|num
14|AC|-1|num
14|AD|0.9750000000000001|num
15|B|6876471|num
15|E|Portimonense|text
15|F|Boavista|text
15|G|1|num
15|H|4|num
15|J|3|num
15|K|A|text
15|L|2.45|num
15|M|3.25|num
15|N|2.875|num
15|O|3|num
15|P|3.3|num
15|Q|2.3|num
15|R|0.25|num
15|S|1.85|num
15|T|2|num
15|U|2.25|num
15|V|1.925|num
15|W|1.925|num
15|X|-1|num
15|Z|1.3|num
15|AA|-1|num
15|AB|1|num
15|AC|0.925|num
15|AD|-1|num
21|F|Benfica|text
23|F|Portimonense|text
31|E|Benfica|text
32|E|Portimonense|text
41|F|Portimonense|text
44|F|Benfica|text
48|B|6875459|num
48|E|Moreirense|text
48|F|SC Farense|text
48|G|1|num
48|H|0|num
48|I|1|num
48|K|H|text
48|L|2.1|num
48|M|3.25|num
48|N|3.5|num
48|O|2.1|num
48|P|3.3|num
48|Q|3.5|num
48|R|-0.25|num
48|S|1.8|num
48|T|2.05|num
48|V|1.95|num
48|W|1.9|num
48|X|1.1|num
48|Y|-1|num
48|AA|0.8|num
48|AB|-1|num
48|AC|-1|num
48|AD|0.8999999999999999|num
49|B|6876499|num
49|E|Estoril|text
49|F|Vizela|text
49|G|2|num
49|H|2|num
49|I|0|num
49|K|D|text
49|L|2.4|num
49|M|3.2|num
49|N|3|num
49|O|2.55|num
49|P|3.2|num
49|Q|2.75|num
49|R|0|num
49|S|1.88|num
49|T|2.02|num
49|V|1.9|num
49|W|1.95|num
49|X|-1|num
49|Y|2.2|num
49|AA|0|num
49|AB|0|num
49|AC|0.8999999999999999|num
49|AD|-1|num
53|E|Benfica|text
53|F|Portimonense|text
57|E|Portimonense|text
58|F|Benfica|text
69|F|Portimonense|text
72|F|Benfica|text
75|E|Benfica|text
77|E|Portimonense|text
84|F|Portimonense|text
86|F|Benfica|text
93|E|Benfica|text
100|E|Portimonense|text
105|F|Benfica|text
107|F|Portimonense|text
111|E|Portimonense|text
112|E|Benfica|text
123|F|Benfica|text
126|F|Portimonense|text
131|E|Portimonense|text
133|B|7513577|num
133|E|Estoril|text
133|F|SC Farense|text
133|G|4|num
133|I|2|num
133|L|2.15|num
133|M|3.6|num
133|N|3.2|num
133|O|1.833|num
133|P|4|num
133|Q|3.8|num
133|R|-0.5|num
133|S|1.875|num
133|T|1.975|num
133|U|2.75|num
133|V|1.875|num
133|W|1.975|num
133|X|0.833|num
133|AA|0.875|num
133|AC|0.875|num
133|AD|-1|num
134|B|7515550|num
134|E|Gil Vicente|text
134|F|Boavista|text
134|G|1|num
134|I|0|num
134|L|2.3|num
134|M|3.5|num
134|N|2.9|num
134|O|2.3|num
134|P|3.3|num
134|Q|3|num
134|R|-0.25|num
134|S|2.05|num
134|T|1.8|num
134|U|2.5|num
134|V|2.05|num
134|W|1.8|num
134|X|1.3|num
134|AA|1.05|num
134|AC|-1|num
134|AD|0.8|num
136|E|Benfica|text
141|F|Portimonense|text
143|F|Benfica|text
146|E|Benfica|text
151|E|Portimonense|text
157|E|Portimonense|text
161|E|Benfica|text
165|F|Benfica|text
167|B|6876607|num
167|E|Arouca|text
167|F|Vizela|text
167|G|5|num
167|I|2|num
167|L|2|num
167|N|3.75|num
167|O|1.8|num
167|P|3.8|num
167|Q|4.2|num
167|R|-0.75|num
167|S|2.05|num
167|T|1.8|num
167|U|2.5|num
167|V|1.85|num
167|W|2|num
167|X|0.8|num
167|AA|1.05|num
167|AC|0.8500000000000001|num
167|AD|-1|num
168|B|6876603|num
168|E|Gil Vicente|text
168|F|Guimaraes|text
168|G|1|num
168|I|0|num
168|L|3.1|num
168|N|2.25|num
168|O|3.25|num
168|P|3.3|num
168|Q|2.25|num
168|R|0.25|num
168|S|1.925|num
168|T|1.925|num
168|U|2.25|num
168|V|1.925|num
168|W|1.925|num
168|X|2.25|num
168|AA|0.925|num
168|AC|-1|num
168|AD|0.925|num
169|F|Portimonense|text
173|E|Benfica|text
177|E|Portimonense|text
181|F|Benfica|text
187|F|Portimonense|text
192|E|Benfica|text
195|E|Portimonense|text
205|E|Portimonense|text
205|F|Benfica|text
213|E|Benfica|text
215|F|Portimonense|text
217|E|Benfica|text
224|E|Portimonense|text
226|F|Benfica|text
233|F|Portimonense|text
236|E|Portimonense|text
243|E|Benfica|text
248|F|Portimonense|text
249|F|Benfica|text
259|E|Benfica|text
260|E|Portimonense|text
267|F|Benfica|text
271|F|Portimonense|text
275|E|Portimonense|text
277|E|Benfica|text
282|B|7096965|num
282|E|Rio Ave|text
282|F|Guimaraes|text
282|G|2|num
282|J|0|num
282|K|H|text
282|L|2.75|num
282|N|2.6|num
282|O|3.3|num
282|P|3.3|num
282|Q|2.2|num
282|R|0.25|num
282|V|2.05|num
282|W|1.8|num
282|X|2.3|num
282|Y|-1|num
282|AA|0.925|num
282|AB|-1|num
282|AC|1.05|num
282|AD|-1|num
283|B|7093688|num
283|E|Boavista|text
283|F|Gil Vicente|text
283|G|1|num
283|J|1|num
283|K|D|text
283|L|2.7|num
283|N|2.625|num
283|O|2.8|num
283|P|3|num
283|Q|2.8|num
283|R|0|num
283|V|2.025|num
283|W|1.825|num
283|X|-1|num
283|Y|2|num
283|AA|0|num
283|AB|0|num
283|AC|-0.5|num
283|AD|0.4125|num
284|F|Benfica|text
289|F|Portimonense|text
291|E|Benfica|text
296|E|Portimonense|text
300|F|Portimonense|text
303|F|Benfica|text
"@

$lines = $changesData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "\|"
    $rowNum = [int]$parts[0]
    $colLetter = $parts[1]
    $rawVal = $parts[2]
    $valType = $parts[3]

    $addr = "$colLetter$rowNum"

    if ($valType -eq "num") {
        $ws.Range($addr).Value = [double]$rawVal
    } else {
        $ws.Range($addr).Value = $rawVal
    }
}

Write-Host "Applied $($lines.Count) cell changes"
